$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert two new list paragraphs right before the paragraph that begins
#    "El actor digita el "Conteo de Efectivo" en el campo. (FA-01)" — i.e.
#    right after the "(EX-01)" paragraph describing the CierreCajaView window.
# ---------------------------------------------------------------------------
$found = $d.Content.Find.Execute("El actor digita el", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate anchor paragraph 'El actor digita el...'"
}
$anchorPara = $d.Content.Find.Parent.Paragraphs.Item(1)
# Walk the Paragraphs collection to find the paragraph object itself.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "El actor digita el*") {
        $targetIndex = $i
        break
    }
}
if ($targetIndex -eq -1) {
    throw "Could not find target paragraph index"
}

$targetPara = $d.Paragraphs.Item($targetIndex)
$insertRange = $targetPara.Range.Duplicate
$insertRange.Collapse(1)

# Insert two blank paragraph marks before the target paragraph; each inherits
# the target paragraph's list formatting (style "Prrafodelista", numId 1).
$insertRange.InsertBefore("`r`r")

$newPara1 = $d.Paragraphs.Item($targetIndex)
$newPara2 = $d.Paragraphs.Item($targetIndex + 1)

$newPara1.Range.Text = "El actor busca en campo " + [char]8220 + "Buscar" + [char]8221 + "." + [char]32 + "(FA-01)"
$newPara2.Range.Text = "El sistema filtra las VENTAs por " + [char]8220 + "No. Venta" + [char]8221 + "."

# ---------------------------------------------------------------------------
# 2) Move the <w:lastRenderedPageBreak/> marker from the "Regresa al flujo
#    normal en el paso 4." paragraph (FA-02, step 3) to the preceding
#    paragraph "El actor hace clic en el botón "Cancelar"." (FA-02, step 2).
# ---------------------------------------------------------------------------
$ns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

$idxClic = -1
$idxRegresa = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "El actor hace clic en el bot*Cancelar*") {
        $idxClic = $i
    }
    if ($t -like "Regresa al flujo normal en el paso 4*") {
        $idxRegresa = $i
    }
}
if ($idxClic -eq -1 -or $idxRegresa -eq -1) {
    throw "Could not locate FA-02 paragraphs"
}

$pClic = $d.Paragraphs.Item($idxClic)
$rngClic = $pClic.Range.Duplicate
$xmlClic = "<w:p $ns><w:pPr><w:pStyle w:val='Prrafodelista'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='15'/></w:numPr><w:jc w:val='both'/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>El actor hace clic en el bot" + [char]0xF3 + "n " + [char]8220 + "Cancelar" + [char]8221 + "</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p>"
$rngClic.InsertXML($xmlClic) | Out-Null

$pRegresa = $d.Paragraphs.Item($idxRegresa)
$rngRegresa = $pRegresa.Range.Duplicate
$xmlRegresa = "<w:p $ns><w:pPr><w:pStyle w:val='Prrafodelista'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='15'/></w:numPr><w:jc w:val='both'/></w:pPr><w:r><w:t>R</w:t></w:r><w:r><w:t>egresa al flujo normal en el</w:t></w:r><w:r><w:t xml:space='preserve'> </w:t></w:r><w:r><w:t>paso</w:t></w:r><w:r><w:t xml:space='preserve'> 4</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p>"
$rngRegresa.InsertXML($xmlRegresa) | Out-Null

Write-Output "Done"
